$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp banner ---
$ws.Range("A1").Value2 = "Datos actualizados a 5 de Mayo de 2020 a las 17:03"

# --- Simple in-place updates (no reordering needed) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value2 = 1215506
$ws.Range("C4").Value2 = 2671
$ws.Range("E4").Value2 = 957305
$ws.Range("G4").Value2 = 211
$ws.Range("H4").Value2 = 70132

# Row 9: Alemania
$ws.Range("B9").Value2 = 166304
$ws.Range("C9").Value2 = 152
$ws.Range("E9").Value2 = 24211

# Row 47: Noruega
$ws.Range("E47").Value2 = 7657
$ws.Range("G47").Value2 = 1
$ws.Range("H47").Value2 = 215

# Row 80: Bulgaria
$ws.Range("B80").Value2 = 1704
$ws.Range("C80").Value2 = 52
$ws.Range("E80").Value2 = 1282
$ws.Range("G80").Value2 = 2
$ws.Range("H80").Value2 = 80

# Row 149: Togo
$ws.Range("B149").Value2 = 128
$ws.Range("C149").Value2 = 2
$ws.Range("E149").Value2 = 45

# --- Somalia update pushes it up above Albania / Sudan / Sri Lanka ---
# Before: 96 Kirguistan, 97 Albania, 98 Sudan, 99 Sri Lanka, 100 Somalia, 101 Niger
# After:  96 Kirguistan, 97 Somalia, 98 Albania, 99 Sudan, 100 Sri Lanka, 101 Niger
$albania = $ws.Range("A97:H97").Value2
$sudan = $ws.Range("A98:H98").Value2
$sriLanka = $ws.Range("A99:H99").Value2

$ws.Range("A97").Value2 = "Somalia"
$ws.Range("B97").Value2 = 835
$ws.Range("C97").Value2 = 79
$ws.Range("D97").Value2 = 75
$ws.Range("E97").Value2 = 722
$ws.Range("F97").Value2 = 2
$ws.Range("G97").Value2 = 3
$ws.Range("H97").Value2 = 38

$ws.Range("A98:H98").Value2 = $albania
$ws.Range("A99:H99").Value2 = $sudan
$ws.Range("A100:H100").Value2 = $sriLanka

# --- Sierra Leona update pushes it up above Islas Feroe / Cabo Verde / Martinica ---
# Before: 135 Congo, 136 Islas Feroe, 137 Cabo Verde, 138 Martinica, 139 Sierra Leona, 140 Liberia
# After:  135 Congo, 136 Sierra Leona, 137 Islas Feroe, 138 Cabo Verde, 139 Martinica, 140 Liberia
$islasFeroe = $ws.Range("A136:H136").Value2
$caboVerde = $ws.Range("A137:H137").Value2
$martinica = $ws.Range("A138:H138").Value2

$ws.Range("A136").Value2 = "Sierra Leona"
$ws.Range("B136").Value2 = 199
$ws.Range("C136").Value2 = 21
$ws.Range("D136").Value2 = 43
$ws.Range("E136").Value2 = 145
$ws.Range("F136").Value2 = 0
$ws.Range("G136").Value2 = 2
$ws.Range("H136").Value2 = 11

$ws.Range("A137:H137").Value2 = $islasFeroe
$ws.Range("A138:H138").Value2 = $caboVerde
$ws.Range("A139:H139").Value2 = $martinica
